# Update column F (dSF) values to reflect repulled/pushed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = -4
    4  = 1
    5  = 1
    6  = -1
    7  = 4
    8  = -3
    9  = -2
    10 = 2
    11 = -1
    12 = -2
    13 = 3
    14 = -6
    15 = 2
    16 = 1
    19 = -1
    21 = 1
    22 = 2
    23 = -1
    27 = -1
    28 = -4
    30 = 2
    32 = 4
    33 = -4
    34 = -1
    35 = 0
    36 = -3
    37 = -2
    38 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
